$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newType = "numeric/character/facotr/Integer"

# Update the "variable type" column. D3/D5 previously held "numeric"; D7:D14
# previously held "Factor". All of these become the new combined label, so
# the underlying shared string for "Factor" is updated in place and reused.
$ws.Range("D3").Value = $newType
$ws.Range("D5").Value = $newType
$ws.Range("D7").Value = $newType
$ws.Range("D8").Value = $newType
$ws.Range("D9").Value = $newType
$ws.Range("D10").Value = $newType
$ws.Range("D11").Value = $newType
$ws.Range("D12").Value = $newType
$ws.Range("D13").Value = $newType
$ws.Range("D14").Value = $newType

# Update the selected/active cell on the sheet
$ws.Range("I14").Select()
